# The workbook gained a new row at the bottom of Sheet1 (B4 = "51847_1"),
# which extends the used range from A1:C3 to A1:C4. The dimension updates
# automatically once the new cell is written.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = "51847_1"

# Mirror the workbook-level "recalc on open" intent from the edit
# (calcPr/@fullCalcOnLoad) by forcing a full recalculation now.
$excel.CalculateFullRebuild()
